$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing F/G values for rows 509-553 (AgTests / AgPosit revisions)
$ws.Range("F509").Value = 9807
$ws.Range("F510").Value = 8012
$ws.Range("F511").Value = 6904
$ws.Range("F512").Value = 8683
$ws.Range("F513").Value = 10584
$ws.Range("F514").Value = 7169
$ws.Range("F515").Value = 5205
$ws.Range("F516").Value = 9533
$ws.Range("F517").Value = 7845
$ws.Range("F518").Value = 7240
$ws.Range("G518").Value = 12
$ws.Range("F519").Value = 8042
$ws.Range("F520").Value = 10464
$ws.Range("F521").Value = 6912
$ws.Range("F522").Value = 5202
$ws.Range("F523").Value = 10300
$ws.Range("F524").Value = 7873
$ws.Range("F525").Value = 7675
$ws.Range("F526").Value = 8844
$ws.Range("F527").Value = 11573
$ws.Range("F528").Value = 8120
$ws.Range("F529").Value = 5732
$ws.Range("F530").Value = 12791
$ws.Range("F531").Value = 9302
$ws.Range("F532").Value = 10263
$ws.Range("F533").Value = 11889
$ws.Range("F534").Value = 16793
$ws.Range("F535").Value = 10135
$ws.Range("G535").Value = 23
$ws.Range("F536").Value = 7988
$ws.Range("F537").Value = 13603
$ws.Range("F538").Value = 11246
$ws.Range("F539").Value = 10524
$ws.Range("F540").Value = 12428
$ws.Range("G540").Value = 64
$ws.Range("F541").Value = 16509
$ws.Range("F542").Value = 10273
$ws.Range("F543").Value = 4665
$ws.Range("F544").Value = 14261
$ws.Range("F545").Value = 16592
$ws.Range("F546").Value = 3784
$ws.Range("F547").Value = 13835
$ws.Range("G547").Value = 151
$ws.Range("F548").Value = 16758
$ws.Range("F549").Value = 10481
$ws.Range("F550").Value = 8415
$ws.Range("F551").Value = 17447
$ws.Range("G551").Value = 196
$ws.Range("F552").Value = 15436
$ws.Range("G552").Value = 172
$ws.Range("F553").Value = 15235
$ws.Range("G553").Value = 177

# Append new rows 554-557 with the latest daily stats
$ws.Range("A554").Value = 44448
$ws.Range("B554").Value = 397759
$ws.Range("C554").Value = 6552
$ws.Range("D554").Value = 377
$ws.Range("E554").Value = 12558
$ws.Range("F554").Value = 16561
$ws.Range("G554").Value = 168

$ws.Range("A555").Value = 44449
$ws.Range("B555").Value = 398278
$ws.Range("C555").Value = 8073
$ws.Range("D555").Value = 519
$ws.Range("E555").Value = 12558
$ws.Range("F555").Value = 18398
$ws.Range("G555").Value = 170

$ws.Range("A556").Value = 44450
$ws.Range("B556").Value = 398690
$ws.Range("C556").Value = 4506
$ws.Range("D556").Value = 412
$ws.Range("E556").Value = 12558
$ws.Range("F556").Value = 9882
$ws.Range("G556").Value = 90

$ws.Range("A557").Value = 44451
$ws.Range("B557").Value = 398744
$ws.Range("C557").Value = 1439
$ws.Range("D557").Value = 54
$ws.Range("E557").Value = 12560
$ws.Range("F557").Value = 6877
$ws.Range("G557").Value = 101

